# Finalize version to deliver reports to user:
# Update the numeric data table (columns B:J, rows 2:43) on the active sheet
# with the refreshed report figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(111,87,102,100,135,174,129,32,81),
    @(4,18,12,7,19,25,22,8,11),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,2,0,2,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(1,0,2,0,2,3,1,1,0),
    @(0,4,4,3,7,13,10,0,4),
    @(3,14,6,4,8,9,9,7,7),
    @(34,28,26,23,44,46,30,7,23),
    @(2,1,2,0,5,4,1,1,1),
    @(7,1,4,1,1,1,6,6,4),
    @(12,9,5,5,2,1,5,0,0),
    @(4,6,7,3,15,14,9,0,6),
    @(8,5,8,10,17,12,5,0,8),
    @(1,6,0,4,4,14,4,0,4),
    @(15,9,22,15,18,27,22,6,15),
    @(5,0,8,4,5,12,10,0,8),
    @(0,0,0,0,0,0,0,0,0),
    @(5,1,7,6,7,6,4,4,3),
    @(4,8,7,5,6,8,7,2,4),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(1,0,0,0,0,1,1,0,0),
    @(30,17,15,35,43,53,32,2,9),
    @(1,1,1,1,2,4,4,0,1),
    @(2,2,0,17,16,17,9,0,0),
    @(10,7,7,5,8,15,8,1,5),
    @(0,0,0,0,1,0,0,0,1),
    @(0,0,0,0,0,0,0,0,0),
    @(17,7,7,12,16,17,11,1,2),
    @(17,3,14,11,5,12,15,7,13),
    @(7,2,4,3,1,5,8,0,4),
    @(10,1,10,8,4,7,7,7,9),
    @(0,0,0,0,0,0,0,0,0),
    @(0,0,0,0,0,0,0,0,0),
    @(11,12,13,9,6,11,8,2,10),
    @(0,0,0,0,0,0,0,0,0),
    @(0,1,5,0,6,5,5,2,3),
    @(3,0,3,0,0,0,2,0,1),
    @(0,0,0,0,0,0,0,0,0),
    @(8,11,5,9,0,6,1,0,6)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $rowIndex = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $colIndex = 2 + $j   # Column B = 2
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$j]
    }
}
